$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert three new columns (C, D, E) - the old column C (event text)
#    shifts right and becomes column F.
# ------------------------------------------------------------------
$ws.Columns("C:E").Insert()

# Match the new columns' width to column B's width (same visual width
# as the other "UN" marker columns).
$ws.Columns("C:E").ColumnWidth = $ws.Columns("B").ColumnWidth

# ------------------------------------------------------------------
# 2) Header row - new weekly columns pushed in front of the old ones.
# ------------------------------------------------------------------
$ws.Range("B1").Value = 'Jun_27'
$ws.Range("C1").Value = 'Jun_26'
$ws.Range("D1").Value = 'Jun_26'
$ws.Range("E1").Value = 'Jun_15'
# F1 already holds 'Jun_10' (shifted automatically by the column insert)

# ------------------------------------------------------------------
# 3) Fill the new C/D/E columns, rows 2-27, with the default 'UN'
#    marker (same as column B), then overwrite the analyst rows that
#    have an actual rating event for the two newest weeks.
# ------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("C$r").Value = 'UN'
    $ws.Range("D$r").Value = 'UN'
    $ws.Range("E$r").Value = 'UN'
}

# Zacks Investment Research (row 5) - Upgrades Hold -> Buy
$ws.Range("C5").Value = '6/26/2018,Upgrades,Hold -> Buy,$17.00'
$ws.Range("D5").Value = '6/26/2018,Upgrades,Hold -> Buy,$17.00'
$ws.Range("C5").Interior.ColorIndex = 35
$ws.Range("D5").Interior.ColorIndex = 35

# Morgan Stanley (row 7) - Raises Target Underweight -> Underweight
$ws.Range("C7").Value = '6/20/2018,Raises Target,Underweight -> Underweight,$8.00 -> $11.00'
$ws.Range("D7").Value = '6/20/2018,Raises Target,Underweight -> Underweight,$8.00 -> $11.00'
$ws.Range("C7").Interior.ColorIndex = 35
$ws.Range("D7").Interior.ColorIndex = 35

# Bank of America (row 20) - Raises Target Buy
$ws.Range("C20").Value = '6/25/2018,Raises Target,Buy,$17.00 -> $20.00'
$ws.Range("D20").Value = '6/25/2018,Raises Target,Buy,$17.00 -> $20.00'
$ws.Range("C20").Interior.ColorIndex = 35
$ws.Range("D20").Interior.ColorIndex = 35

# BidaskClub (row 22) - the previous highlighted "Jun_15" rating now
# lands in the newly inserted E column (the Jun_15 slot); B/C/D revert
# to the plain 'UN' marker with no highlight.
$ws.Range("B22").Value = 'UN'
$ws.Range("B22").Interior.ColorIndex = -4142
$ws.Range("C22").Value = 'UN'
$ws.Range("D22").Value = 'UN'
$ws.Range("E22").Value = '6/13/2018,Upgrades,Buy -> Strong-Buy,'
$ws.Range("E22").Interior.ColorIndex = 35

# ------------------------------------------------------------------
# 4) Two brand-new analyst rows appended at the bottom.
# ------------------------------------------------------------------
$ws.Range("A28").Value = 'Benchmark'
$ws.Range("B28").Value = '6/26/2018,Initiates,Hold,'
$ws.Range("C28").Value = '6/26/2018,Initiates,Hold,'
$ws.Range("D28").Value = '6/26/2018,Initiates,Hold,'

$ws.Range("A29").Value = 'Evercore ISI'
$ws.Range("B29").Value = 'UN'
$ws.Range("C29").Value = 'UN'
$ws.Range("D29").Value = 'UN'
